# Apply updated "dSF" (column F) values to Sheet1 of the workbook.
# These reflect a repull/recalculation of the dSF stat for each outing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -5
    3  = 4
    5  = 6
    6  = 2
    7  = 3
    9  = -2
    10 = 12
    11 = -6
    12 = -1
    13 = 5
    14 = -4
    15 = 2
    16 = -5
    17 = 5
    18 = -2
    19 = -2
    20 = 5
    21 = -2
    22 = -3
    23 = 5
    25 = 4
    26 = 3
    27 = 4
    28 = -2
    29 = 2
    30 = 0
    31 = 3
    32 = 5
    33 = 9
    34 = 2
    35 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
